# Apply Coinranking crypto-price refresh (GitHub Actions bot, 2024-02-29).
# Every write uses a leading apostrophe so the engine stores the literal
# text verbatim (prices like "1.00" / "0.593" / "0.0000136" would
# otherwise be auto-coerced to numbers), then resets .Style to "Normal"
# so the cell keeps the workbook's default (unstyled) cell format, just
# like every other data cell in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'62.650.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.92%  "
$ws.Range("E2").Style = "Normal"

# Row 3: Ethereum
$ws.Range("D3").Value = "'3.461.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.60%  "
$ws.Range("E3").Style = "Normal"

# Row 4: TetherUSD
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"

# Row 5: BNB
$ws.Range("D5").Value = "'411.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.78%  "
$ws.Range("E5").Style = "Normal"

# Row 6: Solana
$ws.Range("D6").Value = "'129.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +15.56%  "
$ws.Range("E6").Style = "Normal"

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "'3.454.21"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.68%  "
$ws.Range("E7").Style = "Normal"

# Row 8: XRP
$ws.Range("D8").Value = "'0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.12%  "
$ws.Range("E8").Style = "Normal"

# Row 9: USDC
$ws.Range("E9").Value = "'  +0.08%  "
$ws.Range("E9").Style = "Normal"

# Row 10: Cardano
$ws.Range("E10").Value = "'  +8.89%  "
$ws.Range("E10").Style = "Normal"

# Row 11: Dogecoin
$ws.Range("E11").Value = "'  +30.48%  "
$ws.Range("E11").Style = "Normal"

# Row 12: Avalanche
$ws.Range("D12").Value = "'43.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +9.23%  "
$ws.Range("E12").Style = "Normal"

# Row 13: TRON
$ws.Range("E13").Value = "'  +0.20%  "
$ws.Range("E13").Style = "Normal"

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'4.008.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.60%  "
$ws.Range("E14").Style = "Normal"

# Row 15: Polkadot
$ws.Range("D15").Value = "'8.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.06%  "
$ws.Range("E15").Style = "Normal"

# Row 16: Chainlink
$ws.Range("D16").Value = "'20.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.91%  "
$ws.Range("E16").Style = "Normal"

# Row 17: WrappedEther
$ws.Range("D17").Value = "'3.367.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.55%  "
$ws.Range("E17").Style = "Normal"

# Row 18: WrappedBTC
$ws.Range("D18").Value = "'62.560.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.08%  "
$ws.Range("E18").Style = "Normal"

# Row 19: Polygon
$ws.Range("E19").Value = "'  +0.73%  "
$ws.Range("E19").Style = "Normal"

# Row 20: Uniswap
$ws.Range("D20").Value = "'10.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.48%  "
$ws.Range("E20").Style = "Normal"

# Row 21: ShibaInu
$ws.Range("D21").Value = "'0.0000136"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +24.56%  "
$ws.Range("E21").Style = "Normal"

# Row 22: ImmutableX
$ws.Range("D22").Value = "'3.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.08%  "
$ws.Range("E22").Style = "Normal"

# Row 23: InternetComputer(DFINITY)
$ws.Range("D23").Value = "'13.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.17%  "
$ws.Range("E23").Style = "Normal"

# Row 24: Litecoin
$ws.Range("D24").Value = "'82.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +9.28%  "
$ws.Range("E24").Style = "Normal"

# Row 25: BitcoinCash
$ws.Range("D25").Value = "'312.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.20%  "
$ws.Range("E25").Style = "Normal"

# Row 26: PancakeSwap
$ws.Range("D26").Value = "'3.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.57%  "
$ws.Range("E26").Style = "Normal"

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'30.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +5.80%  "
$ws.Range("E27").Style = "Normal"

# Row 28: Filecoin
$ws.Range("D28").Value = "'8.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.92%  "
$ws.Range("E28").Style = "Normal"

# Row 29: RenderToken
$ws.Range("E29").Value = "'  +5.99%  "
$ws.Range("E29").Style = "Normal"

# Row 30: Hedera
$ws.Range("E30").Value = "'  +7.73%  "
$ws.Range("E30").Style = "Normal"

# Row 31: Kaspa
$ws.Range("E31").Value = "'  +3.92%  "
$ws.Range("E31").Style = "Normal"

# Row 32: LEO
$ws.Range("E32").Value = "'  -1.55%  "
$ws.Range("E32").Style = "Normal"

# Row 33: InjectiveProtocol
$ws.Range("D33").Value = "'44.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +11.71%  "
$ws.Range("E33").Style = "Normal"

# Row 34: Toncoin
$ws.Range("B34").Value = "'Toncoin"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'2.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +27.24%  "
$ws.Range("E34").Style = "Normal"

# Row 35: Cosmos
$ws.Range("B35").Value = "'Cosmos"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'11.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.48%  "
$ws.Range("E35").Style = "Normal"

# Row 36: Dai
$ws.Range("E36").Value = "'  +0.11%  "
$ws.Range("E36").Style = "Normal"

# Row 37: VeChain
$ws.Range("D37").Value = "'0.0493"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -8.12%  "
$ws.Range("E37").Style = "Normal"

# Row 38: OKB
$ws.Range("D38").Value = "'52.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.19%  "
$ws.Range("E38").Style = "Normal"

# Row 39: LidoDAOToken
$ws.Range("D39").Value = "'3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.57%  "
$ws.Range("E39").Style = "Normal"

# Row 40: FirstDigitalUSD
$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.24%  "
$ws.Range("E40").Style = "Normal"

# Row 41: Stacks
$ws.Range("E41").Value = "'  -5.75%  "
$ws.Range("E41").Style = "Normal"

# Row 42: ARBITRUM
$ws.Range("E42").Value = "'  +4.06%  "
$ws.Range("E42").Style = "Normal"

# Row 43: Celestia
$ws.Range("D43").Value = "'18.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +6.97%  "
$ws.Range("E43").Style = "Normal"

# Row 44: Monero
$ws.Range("B44").Value = "'Monero"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'137.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.05%  "
$ws.Range("E44").Style = "Normal"

# Row 45: Stellar
$ws.Range("B45").Value = "'Stellar"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.125"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.64%  "
$ws.Range("E45").Style = "Normal"

# Row 46: TheGraph
$ws.Range("E46").Value = "'  +3.67%  "
$ws.Range("E46").Style = "Normal"

# Row 47: NEARProtocol
$ws.Range("D47").Value = "'4.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.45%  "
$ws.Range("E47").Style = "Normal"

# Row 48: WEMIXToken
$ws.Range("D48").Value = "'2.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.12%  "
$ws.Range("E48").Style = "Normal"

# Row 49: EnergySwap
$ws.Range("D49").Value = "'22.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.14%  "
$ws.Range("E49").Style = "Normal"

# Row 50: Maker
$ws.Range("D50").Value = "'2.252.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.47%  "
$ws.Range("E50").Style = "Normal"

# Row 51: RocketPoolETH
$ws.Range("D51").Value = "'3.806.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.93%  "
$ws.Range("E51").Style = "Normal"
